$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date for every existing data row (2-386)
# from serial 45202 (2023-10-03) to serial 45203 (2023-10-04).
for ($i = 2; $i -le 386; $i++) {
    $ws.Cells.Item($i, 3).Value2 = 45203
}

# Row 386 gains an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(386).RowHeight = 15

# Append the new record as row 387.
$ws.Range("A387").Value2 = "A 47121-2023"
$ws.Range("B387").Value2 = 45201
$ws.Range("B387").NumberFormat = "YYYY-MM-DD"
$ws.Range("C387").Value2 = 45203
$ws.Range("C387").NumberFormat = "YYYY-MM-DD"
$ws.Range("D387").Value2 = "JÖNKÖPINGS LÄN"
$ws.Range("E387").Value2 = "VAGGERYD"
$ws.Range("G387").Value2 = 7.5
$ws.Range("H387").Value2 = 0
$ws.Range("I387").Value2 = 0
$ws.Range("J387").Value2 = 0
$ws.Range("K387").Value2 = 0
$ws.Range("L387").Value2 = 0
$ws.Range("M387").Value2 = 0
$ws.Range("N387").Value2 = 0
$ws.Range("O387").Value2 = 0
$ws.Range("P387").Value2 = 0
$ws.Range("Q387").Value2 = 0
$ws.Range("R387").WrapText = $true
